$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. Everything currently in columns
# A:W shifts right to B:X (formats, merged cells and the sheet
# dimension all shift automatically).
$ws.Columns("A:A").Insert()

# New column A is the "Match ID" column.
# Row 2 holds the column header (row 1 is the hidden multi-index header,
# row 3 is a hidden spacer row, rows 4-20 are the player data + the
# hidden totals row 20).
$ws.Range("A2").Value2 = "Match ID"
$ws.Range("A4:A20").Value2 = 31

# Match the bold look already used for the rest of the header/index
# column (no border, left as general alignment).
$ws.Range("A2:A20").Font.Bold = $true

# Restore the active selection to the new Match ID column.
$ws.Range("A4:A20").Select()
